# Update MSOSS_title.pptx for the next seminar.
$p = $ppt.ActivePresentation

# --- 1. Refresh the "Updated automatically" date placeholder text that
#        PowerPoint re-stamps on open/save, across the slide master and
#        every custom (slide) layout, from 2022-11-20 to 2023-02-05.
$oldDate = "2022-11-20"
$newDate = "2023-02-05"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- 2. Title slide: swap in the new speaker and talk title.
$s1 = $p.Slides.Item(1)

$speaker = $s1.Shapes.Item("Title 8")
$speaker.TextFrame.TextRange.Text = "Slava Merkin"

$subtitle = $s1.Shapes.Item("Subtitle 11")
$subtitle.TextFrame.TextRange.Text = "Center for Geospace Storms"
# Reset the manual shrink-to-fit scaling that no longer applies to the
# shorter replacement title text.
$subtitle.TextFrame.AutoSize = 2
